$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 196.42857
$ws.Range("I6").Value = 196.42857
$ws.Range("K6").Value = 589.28571
$ws.Range("M6").Value = -477.28571
$ws.Range("H47").Value = 1000
$ws.Range("I47").Value = 1000
$ws.Range("K47").Value = 1000
$ws.Range("M47").Value = -28
$ws.Range("H55").Value = 317.27274
$ws.Range("I55").Value = 298.66666
$ws.Range("J55").Value = 339.6
$ws.Range("K55").Value = 298.66666
$ws.Range("L55").Value = 339.6
$ws.Range("M55").Value = -84.66665999999998
$ws.Range("N55").Value = -767.6
$ws.Range("H80").Value = 7152.143
$ws.Range("I80").Value = 4896
$ws.Range("J80").Value = 8054.6
$ws.Range("K80").Value = 14688
$ws.Range("L80").Value = 24163.8
$ws.Range("M80").Value = -13690
$ws.Range("N80").Value = -26159.8
$ws.Range("H83").Value = 7152.143
$ws.Range("I83").Value = 4896
$ws.Range("J83").Value = 8054.6
$ws.Range("K83").Value = 44064
$ws.Range("L83").Value = 72491.40000000001
$ws.Range("M83").Value = -39072
$ws.Range("N83").Value = -82475.40000000001
$ws.Range("H107").Value = 1705.3334
$ws.Range("I107").Value = 1168.5
$ws.Range("K107").Value = 1168.5
$ws.Range("M107").Value = 751.5
$ws.Range("H125").Value = 3186.4
$ws.Range("J125").Value = 3500
$ws.Range("L125").Value = 31500
$ws.Range("N125").Value = -36420
$ws.Range("H137").Value = 2520.2
$ws.Range("I137").Value = 2400.25
$ws.Range("K137").Value = 7200.75
$ws.Range("M137").Value = -4650.75
$ws.Range("H138").Value = 5174.7856
$ws.Range("J138").Value = 6679.8335
$ws.Range("L138").Value = 20039.5005
$ws.Range("N138").Value = -30319.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2520
$ws.Range("I2").Value = 2687.1428
$ws.Range("J2").Value = 1350
$ws.Range("K2").Value = 2687.1428
$ws.Range("L2").Value = 1350
$ws.Range("M2").Value = -2574.1428
$ws.Range("N2").Value = -1576
$ws.Range("H24").Value = 49838.5
$ws.Range("J24").Value = 49838.5
$ws.Range("L24").Value = 49838.5
$ws.Range("N24").Value = -50586.5
$ws.Range("H100").Value = 49838.5
$ws.Range("J100").Value = 49838.5
$ws.Range("L100").Value = 49838.5
$ws.Range("N100").Value = -52002.5
$ws.Range("H116").Value = 2520
$ws.Range("I116").Value = 2687.1428
$ws.Range("J116").Value = 1350
$ws.Range("K116").Value = 2687.1428
$ws.Range("L116").Value = 1350
$ws.Range("M116").Value = -393.1428000000001
$ws.Range("N116").Value = -5938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2520
$ws.Range("I3").Value = 2687.1428
$ws.Range("J3").Value = 1350
$ws.Range("K3").Value = 2687.1428
$ws.Range("L3").Value = 1350
$ws.Range("M3").Value = -2573.1428
$ws.Range("N3").Value = -1578
$ws.Range("H86").Value = 3400.5557
$ws.Range("I86").Value = 3349.6667
$ws.Range("J86").Value = 3502.3333
$ws.Range("K86").Value = 3349.6667
$ws.Range("L86").Value = 3502.3333
$ws.Range("M86").Value = -2226.6667
$ws.Range("N86").Value = -5748.3333
$ws.Range("H89").Value = 3400.5557
$ws.Range("I89").Value = 3349.6667
$ws.Range("J89").Value = 3502.3333
$ws.Range("K89").Value = 16748.3335
$ws.Range("L89").Value = 17511.6665
$ws.Range("M89").Value = -11132.3335
$ws.Range("N89").Value = -28743.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 15402980
$ws.Range("J4").Value = 21249.092
$ws.Range("L4").Value = 21249.092
$ws.Range("N4").Value = -21473.092
$ws.Range("H7").Value = 1640.7059
$ws.Range("I7").Value = 1162.909
$ws.Range("J7").Value = 2516.6667
$ws.Range("K7").Value = 1162.909
$ws.Range("L7").Value = 2516.6667
$ws.Range("M7").Value = -1049.909
$ws.Range("N7").Value = -2742.6667
$ws.Range("H16").Value = 2834.4614
$ws.Range("I16").Value = 2766.4443
$ws.Range("K16").Value = 2766.4443
$ws.Range("M16").Value = -2479.4443
$ws.Range("H22").Value = 133
$ws.Range("I22").Value = 133
$ws.Range("K22").Value = 133
$ws.Range("M22").Value = 217
$ws.Range("H23").Value = 47240
$ws.Range("J23").Value = 47240
$ws.Range("L23").Value = 47240
$ws.Range("N23").Value = -47720
$ws.Range("H27").Value = 47240
$ws.Range("J27").Value = 47240
$ws.Range("L27").Value = 47240
$ws.Range("N27").Value = -47624
$ws.Range("H113").Value = 2834.4614
$ws.Range("I113").Value = 2766.4443
$ws.Range("K113").Value = 2766.4443
$ws.Range("M113").Value = -596.4443000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 150000
$ws.Range("I22").Value = 150000
$ws.Range("K22").Value = 450000
$ws.Range("M22").Value = -449831
$ws.Range("H23").Value = 511.46155
$ws.Range("I23").Value = 900
$ws.Range("K23").Value = 2700
$ws.Range("M23").Value = -2465
$ws.Range("H27").Value = 150000
$ws.Range("I27").Value = 150000
$ws.Range("K27").Value = 450000
$ws.Range("M27").Value = -449898
$ws.Range("H38").Value = 325.52942
$ws.Range("J38").Value = 255.875
$ws.Range("L38").Value = 767.625
$ws.Range("N38").Value = -1461.625
$ws.Range("H86").Value = 1151
$ws.Range("I86").Value = 1151
$ws.Range("K86").Value = 3453
$ws.Range("M86").Value = -2267
$ws.Range("H89").Value = 1151
$ws.Range("I89").Value = 1151
$ws.Range("K89").Value = 10359
$ws.Range("M89").Value = -4431
$ws.Range("H97").Value = 601.5
$ws.Range("I97").Value = 601.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1804.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1308.5
$ws.Range("N97").ClearContents()
$ws.Range("H109").Value = 2521.4614
$ws.Range("I109").Value = 126.333336
$ws.Range("J109").Value = 3240
$ws.Range("K109").Value = 379.000008
$ws.Range("L109").Value = 9720
$ws.Range("M109").Value = 660.999992
$ws.Range("N109").Value = -11800
$ws.Range("H131").Value = 949.0909
$ws.Range("J131").Value = 991
$ws.Range("L131").Value = 2973
$ws.Range("N131").Value = -13053

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 300000
$ws.Range("I2").Value = 800000
$ws.Range("K2").Value = 800000
$ws.Range("M2").Value = -799888
$ws.Range("H7").Value = 4250
$ws.Range("I7").Value = 4250
$ws.Range("K7").Value = 4250
$ws.Range("M7").Value = -4138
$ws.Range("H12").Value = 25500
$ws.Range("J12").Value = 25500
$ws.Range("L12").Value = 25500
$ws.Range("N12").Value = -25840
$ws.Range("H22").Value = 3899.8
$ws.Range("I22").Value = 1499
$ws.Range("K22").Value = 1499
$ws.Range("M22").Value = -1204
$ws.Range("H27").Value = 3899.8
$ws.Range("I27").Value = 1499
$ws.Range("K27").Value = 1499
$ws.Range("M27").Value = -1392
$ws.Range("H46").Value = 3857.8948
$ws.Range("I46").Value = 3262.5
$ws.Range("K46").Value = 3262.5
$ws.Range("M46").Value = -3074.5
$ws.Range("H61").Value = 1829.2
$ws.Range("J61").Value = 950
$ws.Range("L61").Value = 950
$ws.Range("N61").Value = -1354
$ws.Range("H68").Value = 1800
$ws.Range("I68").Value = 2300
$ws.Range("J68").Value = 1466.6666
$ws.Range("K68").Value = 2300
$ws.Range("L68").Value = 1466.6666
$ws.Range("M68").Value = -1551
$ws.Range("N68").Value = -2964.6666
$ws.Range("H71").Value = 1800
$ws.Range("I71").Value = 2300
$ws.Range("J71").Value = 1466.6666
$ws.Range("K71").Value = 11500
$ws.Range("L71").Value = 7333.333000000001
$ws.Range("M71").Value = -7756
$ws.Range("N71").Value = -14821.333
$ws.Range("H100").Value = 1800
$ws.Range("I100").Value = 1800
$ws.Range("K100").Value = 1800
$ws.Range("M100").Value = -1259
$ws.Range("H113").Value = 1829.2
$ws.Range("J113").Value = 950
$ws.Range("L113").Value = 950
$ws.Range("N113").Value = -5290
$ws.Range("H122").Value = 3733
$ws.Range("I122").Value = 3733
$ws.Range("K122").Value = 11199
$ws.Range("M122").Value = -8749
$ws.Range("H126").Value = 4250
$ws.Range("I126").Value = 4250
$ws.Range("K126").Value = 12750
$ws.Range("M126").Value = -10280
